# Add a new "TP5" evaluation row to the grade sheet (Sheet1), matching the
# formatting of the existing TP rows, and refresh the selection that Excel
# had recorded at save time.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 previously held just a lone space in B9 (a placeholder). Turn it into
# a full TP entry like rows 5/7 (TP1/TP3): copy their cell formatting into
# row 9's A/B/C cells so the new row matches the sheet's existing look
# (fill colour, alignment, wrap text, number format, etc.).
$ws.Range("A7").Copy()
$ws.Range("A9").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B7").Copy()
$ws.Range("B9").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("C7").Copy()
$ws.Range("C9").PasteSpecial(-4122)   # xlPasteFormats

# Fill in the new TP5 data.
$ws.Range("A9").Value = "TP5"
$ws.Range("B9").Value = "Identificar erro no padrão MTV (os dois apresentaram)`n- Framework Jungle para Python`n- Verifica falhas na implementação`n- Ponto positivo: focaram em um padrão que não vimos em sala. Apontou limitações.`n- Ponto negativo: a explicação do que foi feito ficou fraca, pois mais embassamento do MTV e talvez um toy example de uma aplicação simples devia ter sido dada. Pois ainda tinham 2 minutos."
$ws.Range("C9").Value = 7

# Row grew tall to fit the new wrapped comment text.
$ws.Rows.Item(9).RowHeight = 170

# TOTAL (C2 = SUM(C4:C999)) recalculates automatically once C9 has a value.

# Author's selection moved to C10 (just below the new row) when the file
# was last saved.
$ws.Range("C10").Select() | Out-Null
